$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = "123456789"
$ws.Range("B4").Value = "a"
$ws.Range("C4").Value = "b"
$ws.Range("D4").Value = "1"
$ws.Range("E4").Value = $false

# Row 5
$ws.Range("A5").Value = "987654324"
$ws.Range("B5").Value = "ss"
$ws.Range("C5").Value = "sss"
$ws.Range("D5").Value = "ss"
$ws.Range("E5").Value = $false

# Row 6
$ws.Range("A6").Value = "123333223"
$ws.Range("B6").Value = "aa"
$ws.Range("C6").Value = "aaaaaa"
$ws.Range("D6").Value = "aa"
$ws.Range("E6").Value = $false

# Row 7
$ws.Range("A7").Value = "316440262"
$ws.Range("B7").Value = "shani"
$ws.Range("C7").Value = "waizman"
$ws.Range("D7").Value = "0393"
$ws.Range("E7").Value = $true
